$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.486451983451843
$ws.Range("B1").Value = 2.85660195350647
$ws.Range("D1").Value = 1.366375923156738
$ws.Range("E1").Value = 0.8627302050590515
